$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset (Albahaca / Feria Lagunitas de Puerto Montt) gains two new
# weekly price rows. Insert the first new row above the current row 59
# (pushing rows 59-85 down to 60-86), then insert the second new row above
# the current row 87 (pushing the former rows 86-87 down to 88-89).

$ws.Rows.Item(59).Insert()
$ws.Rows.Item(87).Insert()

# --- New row 59 ---
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 44567
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112052
$ws.Cells.Item(59, 7).Value = "Albahaca"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 80
$ws.Cells.Item(59, 11).Value = 7000
$ws.Cells.Item(59, 12).Value = 7000
$ws.Cells.Item(59, 13).Value = 7000
$ws.Cells.Item(59, 14).Value = "$/docena de matas"
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 1167
$ws.Cells.Item(59, 17).Value = 6
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# --- New row 87 ---
$ws.Cells.Item(87, 1).Value = 4
$ws.Cells.Item(87, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87, 3).Value = "Los Lagos"
$ws.Cells.Item(87, 4).Value = 44568
$ws.Cells.Item(87, 5).Value = 10
$ws.Cells.Item(87, 6).Value = 100112052
$ws.Cells.Item(87, 7).Value = "Albahaca"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 80
$ws.Cells.Item(87, 11).Value = 9000
$ws.Cells.Item(87, 12).Value = 9000
$ws.Cells.Item(87, 13).Value = 9000
$ws.Cells.Item(87, 14).Value = "$/docena de matas"
$ws.Cells.Item(87, 15).Value = "Región Metropolitana"
$ws.Cells.Item(87, 16).Value = 1500
$ws.Cells.Item(87, 17).Value = 6
$ws.Cells.Item(87, 18).Value = "Hortaliza"
